# Update the NATMI LR-pair sheet (Fgf2-Fgfr2) with newly recomputed TPM-based
# statistics. Only the numeric score columns (G..T) change; the identifying
# columns (A..F) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> hashtable of column letter -> new value.
$updates = @{
    2  = @{ G=0.1636683333333333; H=0.491005; I=0.008639493057305454; J=0.008639493057305455;
            K=2; L=0.6666666666666666; M=0.05601; N=0.16803;
            O=0.02710547761971223; P=0.02710547761971223;
            Q=0.009167063350000001; R=0.08250357015000001;
            S=0.0002341775857104522; T=0.0002341775857104522 }

    3  = @{ G=0.1636683333333333; H=0.491005; I=0.008639493057305454; J=0.008639493057305455;
            N=5.594253;
            O=0.902427539668559; P=0.9024275396685592;
            Q=0.3052006882516667; R=2.746806194265;
            S=0.007796516463687757; T=0.00779651646368776 }

    4  = @{ G=0.1636683333333333; H=0.491005; I=0.008639493057305454; J=0.008639493057305455;
            M=0.145611; N=0.436833;
            O=0.07046698271172858; P=0.07046698271172858;
            Q=0.023831909685; R=0.214487187165;
            S=0.0006087990079072425; T=0.0006087990079072426 }

    5  = @{ I=0.808839719627903; J=0.8088397196279031;
            K=2; L=0.6666666666666666; M=0.05601; N=0.16803;
            O=0.02710547761971223; P=0.02710547761971223;
            Q=0.8582314842600001; R=7.724083358340001;
            S=0.02192398691830844; T=0.02192398691830845 }

    6  = @{ I=0.808839719627903; J=0.8088397196279031;
            N=5.594253;
            O=0.902427539668559; P=0.9024275396685592;
            Q=28.573255106326; R=257.159295956934;
            S=0.7299192381700156; T=0.7299192381700158 }

    7  = @{ I=0.808839719627903; J=0.8088397196279031;
            M=0.145611; N=0.436833;
            O=0.07046698271172858; P=0.07046698271172858;
            Q=2.231172016686001; R=20.080548150174;
            S=0.05699649453957883; T=0.05699649453957884 }

    8  = @{ G=3.457711333333334; H=10.373134; I=0.1825207873147914; J=0.1825207873147914;
            K=2; L=0.6666666666666666; M=0.05601; N=0.16803;
            O=0.02710547761971223; P=0.02710547761971223;
            Q=0.19366641178; R=1.74299770602;
            S=0.004947313115693335; T=0.004947313115693335 }

    9  = @{ G=3.457711333333334; H=10.373134; I=0.1825207873147914; J=0.1825207873147914;
            N=5.594253;
            O=0.902427539668559; P=0.9024275396685592;
            Q=6.447770666544668; R=58.029935998902;
            S=0.1647117850348556; T=0.1647117850348556 }

    10 = @{ G=3.457711333333334; H=10.373134; I=0.1825207873147914; J=0.1825207873147914;
            M=0.145611; N=0.436833;
            O=0.07046698271172858; P=0.07046698271172858;
            Q=0.5034808049580001; R=4.531327244622;
            S=0.0128616891642425; T=0.0128616891642425 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
